$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the unused crimp-contact row (old row 20: Harwin FEMALE CRIMP CONTACT TIN LOOSE)
$ws.Rows("20:20").Delete()

# The delete operation leaves a dangling #REF! in the shared "A" running-number
# formula column because the cell that used to anchor row 21's formula was the
# deleted row. Re-apply the simple running total formula so it again reads
# "previous row + 1" like the rest of the column.
$ws.Range("A20").Formula = "=A19+1"

# Move the active selection to A20, matching the post-edit cursor position.
$ws.Range("A20").Select()
